# Edit script: apply the changes described by the commit diff.
#
# Summary of the semantic edit:
#  - Delete the slide "Is the volume of movies released related to movie
#    rating?" (display position 9). All later slides shift up by one, and
#    PowerPoint's own bookkeeping (notesMaster rel id, sldIdLst r:ids,
#    notes-slide text rotation, creationIds, etc.) follows automatically.
#  - On the "Motivation & Summary" slide, reword one bullet from
#    "Volume of movies released in a year" to "Year".
#  - On the "Data Analysis" slide (Question 1 bullet) and the dedicated
#    "Question 1" slide, reword
#    "Is the volume of movies released related to movie rating?" to
#    "Is the movie release year related to movie rating?".

$p = $ppt.ActivePresentation

# --- 1) Delete the "Is the volume of movies released..." slide (slide 9) ---
$s9 = $p.Slides.Item(9)
$s9.Delete()

# --- 2) "Motivation & Summary" slide (slide 2): reword bullet ---
$s2 = $p.Slides.Item(2)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $sh = $s2.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -like "*Volume of movies released in a year*") {
            $tr.Replace("Volume of movies released in a year", "Year") | Out-Null
        }
    }
}

# --- 3) "Data Analysis" slide (slide 7): reword Question 1 bullet ---
$s7 = $p.Slides.Item(7)
for ($i = 1; $i -le $s7.Shapes.Count; $i++) {
    $sh = $s7.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -like "*Is the volume of movies released related to movie rating?*") {
            $tr.Replace("Is the volume of movies released related to movie rating?", "Is the movie release year related to movie rating?") | Out-Null
        }
    }
}

# --- 4) "Question 1" slide (slide 8): reword subtitle ---
$s8 = $p.Slides.Item(8)
for ($i = 1; $i -le $s8.Shapes.Count; $i++) {
    $sh = $s8.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -like "*Is the volume of movies released related to movie rating?*") {
            $tr.Replace("Is the volume of movies released related to movie rating?", "Is the movie release year related to movie rating?") | Out-Null
        }
    }
}
